$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp header string ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 07:50"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 245373
$ws.Range("C4").Value = 496
$ws.Range("E4").Value = 228875

# --- Row 35: Japon ---
$ws.Range("C35").Value = 0
$ws.Range("G35").Value = 0

# --- Row 38: Pakistan ---
$ws.Range("B38").Value = 2441
$ws.Range("C38").Value = 20
$ws.Range("E38").Value = 2281
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 35

# --- Rows 74/75: Bulgaria overtakes Letonia in ranking (countries swap rows) ---
# Row 74 becomes Bulgaria, with updated figures
$ws.Range("A74").Value = "Bulgaria"
$ws.Range("B74").Value = 477
$ws.Range("C74").Value = 20
$ws.Range("D74").Value = 30
$ws.Range("E74").Value = 435
$ws.Range("F74").Value = 18
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 12

# Row 75 becomes Letonia, keeping its (unchanged) prior figures
$ws.Range("A75").Value = "Letonia"
$ws.Range("B75").Value = 458
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 31
$ws.Range("E75").Value = 427
$ws.Range("F75").Value = 3
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0

# --- Row 124: Paraguay ---
$ws.Range("D124").Value = 4
$ws.Range("E124").Value = 85
